# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the data rows (2-11): each row keeps its
# static descriptive columns (A,B,C,E,F,G,H,I,J,K) but the per-record
# columns (D,L,M,N,O,P,Q,R,S,T) are re-assigned according to the mapping
# below (new row <- old row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> old row number
$map = @{
    2  = 6
    3  = 9
    4  = 7
    5  = 4
    6  = 8
    7  = 10
    8  = 11
    9  = 3
    10 = 2
    11 = 5
}

# Columns whose values move together with the row mapping.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current ("before") values for every relevant column/row
# before any writes happen, since the mapping's source and destination
# rows overlap.
$snapshot = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Apply the permutation: new row's values <- snapshot of old row.
foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $src[$col]
    }
}
